$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '25.833.32'
$ws.Cells.Item(2, 5).Value = '  -0.56%  '
$ws.Cells.Item(3, 5).Value = '  -0.60%  '
$ws.Cells.Item(4, 5).Value = '  +0.36%  '
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = '@'
$c.Value = '215.41'
$c.Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  +0.10%  '
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = '@'
$c.Value = '0.5103'
$c.Style = 'Normal'
$ws.Cells.Item(6, 5).Value = '  -0.56%  '
$ws.Cells.Item(7, 5).Value = '  +0.35%  '
$ws.Cells.Item(8, 5).Value = '  +0.42%  '
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = '@'
$c.Value = '0.06395'
$c.Style = 'Normal'
$ws.Cells.Item(9, 5).Value = '  +0.68%  '
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = '@'
$c.Value = '19.39'
$c.Style = 'Normal'
$ws.Cells.Item(10, 5).Value = '  -2.06%  '
$ws.Cells.Item(11, 5).Value = '  +0.21%  '
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = '@'
$c.Value = '4.266'
$c.Style = 'Normal'
$ws.Cells.Item(12, 5).Value = '  -0.25%  '
$ws.Cells.Item(13, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(13, 4).Value = '1.853.63'
$ws.Cells.Item(13, 5).Value = '  -0.68%  '
$ws.Cells.Item(14, 2).Value = 'WrappedEther'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(14, 4).Value = '1.625.36'
$ws.Cells.Item(14, 5).Value = '  -0.79%  '
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = '@'
$c.Value = '0.5600'
$c.Style = 'Normal'
$ws.Cells.Item(15, 5).Value = '  +2.51%  '
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = '@'
$c.Value = '63.41'
$c.Style = 'Normal'
$ws.Cells.Item(16, 5).Value = '  -1.46%  '
$ws.Cells.Item(17, 4).Value = '0.0₅7545'
$ws.Cells.Item(17, 5).Value = '  -2.78%  '
$ws.Cells.Item(18, 4).Value = '25.835.14'
$ws.Cells.Item(18, 5).Value = '  -0.60%  '
$ws.Cells.Item(19, 5).Value = '  +0.68%  '
$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = '@'
$c.Value = '194.37'
$c.Style = 'Normal'
$ws.Cells.Item(20, 5).Value = '  -1.68%  '
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = '@'
$c.Value = '4.329'
$c.Style = 'Normal'
$ws.Cells.Item(21, 5).Value = '  -2.62%  '
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = '@'
$c.Value = '9.804'
$c.Style = 'Normal'
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = '@'
$c.Value = '6.001'
$c.Style = 'Normal'
$ws.Cells.Item(23, 5).Value = '  -1.33%  '
$ws.Cells.Item(24, 5).Value = '  +0.27%  '
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = '@'
$c.Value = '1.824'
$c.Style = 'Normal'
$ws.Cells.Item(25, 5).Value = '  -5.55%  '
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = '@'
$c.Value = '0.1286'
$c.Style = 'Normal'
$ws.Cells.Item(26, 5).Value = '  +4.31%  '
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = '@'
$c.Value = '141.21'
$c.Style = 'Normal'
$ws.Cells.Item(27, 5).Value = '  -0.19%  '
$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = '@'
$c.Value = '6.751'
$c.Style = 'Normal'
$ws.Cells.Item(28, 5).Value = '  -1.42%  '
$ws.Cells.Item(29, 5).Value = '  -1.24%  '
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = '@'
$c.Value = '1.237'
$c.Style = 'Normal'
$ws.Cells.Item(30, 5).Value = '  -0.30%  '
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = '@'
$c.Value = '0.04892'
$c.Style = 'Normal'
$ws.Cells.Item(31, 5).Value = '  +0.95%  '
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = '@'
$c.Value = '3.309'
$c.Style = 'Normal'
$ws.Cells.Item(32, 5).Value = '  +0.74%  '
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = '@'
$c.Value = '3.190'
$c.Style = 'Normal'
$ws.Cells.Item(33, 5).Value = '  -0.88%  '
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = '@'
$c.Value = '1.557'
$c.Style = 'Normal'
$ws.Cells.Item(34, 5).Value = '  +1.18%  '
$ws.Cells.Item(35, 5).Value = '  +0.24%  '
$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = '@'
$c.Value = '0.8955'
$c.Style = 'Normal'
$ws.Cells.Item(36, 5).Value = '  -2.16%  '
$ws.Cells.Item(37, 4).Value = '1.134.06'
$ws.Cells.Item(37, 5).Value = '  +0.94%  '
$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = '@'
$c.Value = '2.544'
$c.Style = 'Normal'
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = '@'
$c.Value = '0.5484'
$c.Style = 'Normal'
$ws.Cells.Item(39, 5).Value = '  -1.36%  '
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = '@'
$c.Value = '0.01561'
$c.Style = 'Normal'
$ws.Cells.Item(40, 5).Value = '  -0.67%  '
$ws.Cells.Item(41, 5).Value = '  +0.33%  '
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = '@'
$c.Value = '5.586'
$c.Style = 'Normal'
$ws.Cells.Item(42, 5).Value = '  +0.28%  '
$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = '@'
$c.Value = '0.7961'
$c.Style = 'Normal'
$ws.Cells.Item(43, 5).Value = '  -1.28%  '
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = '@'
$c.Value = '97.30'
$c.Style = 'Normal'
$ws.Cells.Item(44, 5).Value = '  -2.05%  '
$ws.Cells.Item(45, 4).Value = '1.777.13'
$ws.Cells.Item(45, 5).Value = '  -0.09%  '
$ws.Cells.Item(46, 5).Value = '  -7.26%  '
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = '@'
$c.Value = '0.4432'
$c.Style = 'Normal'
$ws.Cells.Item(47, 5).Value = '  -2.23%  '
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = '@'
$c.Value = '54.91'
$c.Style = 'Normal'
$ws.Cells.Item(48, 5).Value = '  -0.34%  '
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = '@'
$c.Value = '7.579'
$c.Style = 'Normal'
$ws.Cells.Item(50, 5).Value = '  +1.42%  '
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = '@'
$c.Value = '1.002'
$c.Style = 'Normal'
$ws.Cells.Item(51, 5).Value = '  +0.00%  '
